$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume table in place. Price strings that look
# like plain decimals (e.g. "247.98") would otherwise be auto-converted to
# numbers by Excel's normal text-to-number inference, so those are written
# with a leading apostrophe to force literal text, then the cell style is
# reset back to "Normal" so no stray quote-prefix formatting sticks around.
$ws.Range("D2").Value = '30.826.92'
$ws.Range("E2").Value = '  +1.36%  '
$ws.Range("D3").Value = '1.892.96'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''247.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '''0.4742'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''0.2940'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.80%  '
$ws.Range("D9").Value = '''0.06545'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("D10").Value = '''22.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.06%  '
$ws.Range("D11").Value = '''0.07800'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.06%  '
$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").Value = '''97.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '''0.7422'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.56%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.892.09'
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").Value = '''5.272'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.64%  '
$ws.Range("D16").Value = '''287.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.63%  '
$ws.Range("D17").Value = '30.828.03'
$ws.Range("E17").Value = '  +1.41%  '
$ws.Range("D18").Value = '''13.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("D19").Value = '''0.000007548'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").Value = '''1.0000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '2.141.71'
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("D22").Value = '''5.355'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.57%  '
$ws.Range("D23").Value = '''1.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '''6.287'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.93%  '
$ws.Range("D25").Value = '''9.252'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("D26").Value = '''164.85'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("D27").Value = '''19.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.61%  '
$ws.Range("D28").Value = '''1.927'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("D29").Value = '''1.344'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("D30").Value = '''0.09796'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.27%  '
$ws.Range("D31").Value = '''1.493'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("D32").Value = '''4.332'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.30%  '
$ws.Range("D33").Value = '''4.203'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("D34").Value = '''0.04914'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.86%  '
$ws.Range("D35").Value = '''1.133'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.28%  '
$ws.Range("D36").Value = '''0.7021'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.27%  '
$ws.Range("D37").Value = '''2.726'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D38").Value = '''0.01905'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.09%  '
$ws.Range("D39").Value = '''2.816'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.75%  '
$ws.Range("D40").Value = '''76.89'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.81%  '
$ws.Range("D41").Value = '''6.343'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.61%  '
$ws.Range("D42").Value = '''2.009'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.28%  '
$ws.Range("D43").Value = '''0.4302'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.76%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '''0.8413'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '''1.001'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = '''102.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("D47").Value = '''9.649'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.59%  '
$ws.Range("D48").Value = '''7.067'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("D49").Value = '''35.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("D50").Value = '''914.85'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '''0.3984'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.05%  '
